$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force the cell to keep / become a text (string) value even when the
    # text looks like a number (e.g. "0", "130.0000"), while preserving the
    # cell's original number format / style.
    $fmt = $range.NumberFormat()
    $range.NumberFormat = "@"
    $range.Value() = $text
    $range.NumberFormat = $fmt
}

# --- Update row 7 (item #1): AMARYL 1MG 30 TAB -> ORGASOL LIGHT CREAM ---
$ws.Range("C7").Value() = "ORGASOL LIGHT CREAM"
$ws.Range("H7").Value() = "0:0"
Set-TextValue $ws.Range("L7") "0"
$ws.Range("N7").Value() = "130.00"
Set-TextValue $ws.Range("P7") "130.0000"
$ws.Range("Q7").Value() = "1:0"

# --- Update row 8 (item #2): COLONA 30 F.C.TAB -> PRISBRINA  CAPS ---
$ws.Range("C8").Value() = "PRISBRINA  CAPS"
$ws.Range("H8").Value() = "0:0"
Set-TextValue $ws.Range("L8") "0"
$ws.Range("N8").Value() = "150.00"
Set-TextValue $ws.Range("P8") "150.0000"
$ws.Range("Q8").Value() = "1:0"

# --- Update row 9 (item #3): CONTROLOC ... -> QUICK NAIL  LOTION ---
$ws.Range("C9").Value() = "QUICK NAIL  LOTION"
$ws.Range("H9").Value() = "0:0"
Set-TextValue $ws.Range("L9") "0"
$ws.Range("N9").Value() = "85.00"
Set-TextValue $ws.Range("P9") "1955.0000"
$ws.Range("Q9").Value() = "23:0"

# --- Remove old rows 10-14 (items #4-#8), shifting the totals/footer rows up ---
$ws.Range("A10:Q14").EntireRow.Delete()

# --- Update the totals row (now row 10) ---
$ws.Range("P10").Value() = 2235

# --- Update the footer timestamp (now row 11) ---
$ws.Range("A11").Value() = "Saturday, 24 May, 2025 10:18 AM"

$wb.Save()
